# Generate Report for Handoff
# Updates the localization-status report from the 94c0e964... source file
# to the 232aaed9... source file: refreshes the UUID-based file names, the
# xliff handoff filenames/hashes, the handoff/handback timestamps, clears
# the (now stale) "Latest Target File" / "Latest Handback File" columns,
# and drops the now-broken "Latest Target File" hyperlinks on the
# language sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "94c0e964-2323-4f37-ae93-3328c6e77f30"
$newGuid = "232aaed9-32b6-4121-88c8-e94874ec313b"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 20:59:57"

# Hyperlinks can only be refreshed by re-adding them (editing
# TextToDisplay in place leaves the old link behind), so recreate the
# single B2 link, keeping its original target address.
$overviewTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41575a505568da6c8733ba8bbc79ad2a9a161c22/e2e/$oldGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewTarget, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.17cd0587658399d6fa67e6a95b1b145583560315.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-22 20:59:52"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Columns.Item(9).ColumnWidth = 17.8
$wsZh.Columns.Item(10).ColumnWidth = 20.8

# Drop the now-dead "Latest Target File" (I2) hyperlink; keep A2's link
# (re-added with refreshed display text / same original target).
$zhTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41575a505568da6c8733ba8bbc79ad2a9a161c22/e2e/$oldGuid.md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhTarget, "", "", "$newGuid.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.17cd0587658399d6fa67e6a95b1b145583560315.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-22 20:59:57"
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Columns.Item(9).ColumnWidth = 17.8
$wsDe.Columns.Item(10).ColumnWidth = 20.8

$deTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41575a505568da6c8733ba8bbc79ad2a9a161c22/e2e/$oldGuid.md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deTarget, "", "", "$newGuid.md")
